$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date for eb9e7843-...md (row 4)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-25 16:48:12"

# "zh-cn" sheet: Correspond Handoff / Handback datetimes for the
# eb9e7843-... xlf row (row 4)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-25 16:48:02"
$wsZhCn.Range("K4").Value = "2016-08-25 16:48:31"

# "de-de" sheet: Correspond Handoff datetime (shared text with the
# Overview sheet's Latest HO Xliff Generate Date) and Correspond
# Handback datetime for the eb9e7843-... xlf row (row 4)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-25 16:48:12"
$wsDeDe.Range("K4").Value = "2016-08-25 16:48:38"
